# Updates applied to cryptos worksheet on Sat Sep 28 08:30:17 UTC 2024 (prices + 1h volume %)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "65.668.66"
    "E2" = "  +0.01%  "
    "D3" = "2.671.90"
    "E3" = "  -0.01%  "
    "E4" = "  +0.07%  "
    "D5" = "600.00"
    "E5" = "  -1.21%  "
    "D6" = "155.68"
    "E6" = "  -1.80%  "
    "E7" = "  +0.04%  "
    "D8" = "0.603"
    "E8" = "  +2.00%  "
    "E9" = "  -2.35%  "
    "D10" = "5.90"
    "E10" = "  +1.23%  "
    "E11" = "  -3.15%  "
    "D12" = "0.154"
    "E12" = "  -0.14%  "
    "D13" = "29.32"
    "E13" = "  -1.72%  "
    "D14" = "0.0000195"
    "E14" = "  -0.23%  "
    "D15" = "3.152.24"
    "E15" = "  -0.01%  "
    "D16" = "65.505.46"
    "E16" = "  +0.18%  "
    "D17" = "2.670.65"
    "E17" = "  -0.07%  "
    "D18" = "12.53"
    "E18" = "  -2.33%  "
    "D19" = "4.80"
    "E19" = "  -2.33%  "
    "D20" = "7.51"
    "E20" = "  +1.62%  "
    "D21" = "350.07"
    "E21" = "  -3.18%  "
    "E22" = "  -0.12%  "
    "D23" = "70.19"
    "E23" = "  +1.62%  "
    "D24" = "9.77"
    "E24" = "  +1.74%  "
    "D25" = "0.0000109"
    "E25" = "  +1.69%  "
    "D26" = "1.64"
    "E26" = "  -3.66%  "
    "E27" = "  -2.24%  "
    "E28" = "  +1.48%  "
    "D29" = "8.09"
    "E29" = "  -2.23%  "
    "E30" = "  +0.10%  "
    "D31" = "536.59"
    "E31" = "  -0.98%  "
    "D32" = "2.15"
    "E32" = "  -3.22%  "
    "E33" = "  -5.00%  "
    "E34" = "  +2.07%  "
    "D35" = "5.40"
    "E35" = "  -4.49%  "
    "D36" = "0.422"
    "E36" = "  -2.95%  "
    "D37" = "20.35"
    "E37" = "  -1.71%  "
    "D38" = "159.54"
    "E38" = "  -2.03%  "
    "D39" = "0.999"
    "E39" = "  +0.01%  "
    "D40" = "1.94"
    "E40" = "  -4.53%  "
    "D41" = "1.00"
    "E41" = "  +0.09%  "
    "D42" = "42.42"
    "E42" = "  -0.06%  "
    "D43" = "165.71"
    "E43" = "  -1.21%  "
    "D44" = "4.08"
    "E44" = "  -2.98%  "
    "D45" = "0.0613"
    "E45" = "  -0.89%  "
    "D46" = "23.00"
    "E46" = "  -1.11%  "
    "D47" = "2.22"
    "E47" = "  -5.28%  "
    "E48" = "  -2.48%  "
    "D49" = "0.0260"
    "E49" = "  -2.51%  "
    "D50" = "0.0995"
    "E50" = "  +0.30%  "
    "D51" = "20.01"
    "E51" = "  +0.46%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "600.00", "1.00")
    # keep their exact original formatting instead of being coerced to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}
